$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.400.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.19%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.868.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.51%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.7033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.07963"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.45%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.3133"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'24.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.90%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.935.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'93.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.13%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.169"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.12%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.7021"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.46%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.489"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.02%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.000008657"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.88%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'29.461.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.43%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'251.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.38%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.143.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.64%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -1.27%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.03%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.648"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.13%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -2.87%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.000"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.42%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'161.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.46%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'18.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.21%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.310"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.259"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.215"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.23%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05266"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.82%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.905"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.32%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.7576"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.52%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D37").Value = "'2.707"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.04%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.280.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.39%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01875"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.772"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.8952"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.67%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'109.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.81%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'6.019"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.67%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'71.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.17%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'RocketPoolETH"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.042.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.69%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'BabyDogeCoin"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000127"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.59%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.597"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.94%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.5183"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.88%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.4293"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.16%  "
$ws.Range("E51").Style = "Normal"
